$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.768.07'
$ws.Cells.Item(2, 5).Value = '  -2.61%  '
$ws.Cells.Item(3, 4).Value = '1.781.46'
$ws.Cells.Item(3, 5).Value = '  -2.10%  '
$ws.Cells.Item(4, 5).Value = '  +0.18%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '310.97'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  +0.11%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5115'
$ws.Cells.Item(7, 5).Value = '  -0.82%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3767'
$ws.Cells.Item(8, 5).Value = '  -2.73%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07771'
$ws.Cells.Item(9, 5).Value = '  -7.98%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.17'
$ws.Cells.Item(10, 5).Value = '  -1.75%  '
$ws.Cells.Item(11, 5).Value = '  -2.30%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.201'
$ws.Cells.Item(13, 5).Value = '  -3.16%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.11'
$ws.Cells.Item(14, 5).Value = '  -4.25%  '
$ws.Cells.Item(15, 4).Value = '1.778.87'
$ws.Cells.Item(15, 5).Value = '  -1.85%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '7.164'
$ws.Cells.Item(16, 5).Value = '  -4.44%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '91.89'
$ws.Cells.Item(17, 5).Value = '  -0.84%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001070'
$ws.Cells.Item(18, 5).Value = '  -5.93%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06544'
$ws.Cells.Item(19, 5).Value = '  -2.28%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.003'
$ws.Cells.Item(20, 5).Value = '  +0.15%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '16.94'
$ws.Cells.Item(21, 5).Value = '  -4.51%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.917'
$ws.Cells.Item(22, 5).Value = '  -2.65%  '
$ws.Cells.Item(23, 4).Value = '27.821.75'
$ws.Cells.Item(23, 5).Value = '  -2.49%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '10.95'
$ws.Cells.Item(24, 5).Value = '  -4.29%  '
$ws.Cells.Item(25, 5).Value = '  -1.23%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '158.62'
$ws.Cells.Item(26, 5).Value = '  +0.18%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.17'
$ws.Cells.Item(27, 5).Value = '  -4.46%  '
$ws.Cells.Item(28, 4).Value = '1.985.10'
$ws.Cells.Item(28, 5).Value = '  -1.92%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.346'
$ws.Cells.Item(29, 5).Value = '  -3.21%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '123.48'
$ws.Cells.Item(30, 5).Value = '  -2.04%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.1075'
$ws.Cells.Item(31, 5).Value = '  -0.53%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.030'
$ws.Cells.Item(32, 5).Value = '  -5.90%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.614'
$ws.Cells.Item(33, 5).Value = '  -1.96%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.474'
$ws.Cells.Item(34, 5).Value = '  -4.83%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.07038'
$ws.Cells.Item(35, 5).Value = '  -5.75%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02299'
$ws.Cells.Item(36, 5).Value = '  -2.76%  '
$ws.Cells.Item(37, 2).Value = 'FraxShare'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '8.688'
$ws.Cells.Item(37, 5).Value = '  -0.72%  '
$ws.Cells.Item(38, 2).Value = 'Algorand'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2116'
$ws.Cells.Item(38, 5).Value = '  -4.94%  '
$ws.Cells.Item(39, 2).Value = 'Aptos'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '11.54'
$ws.Cells.Item(39, 5).Value = '  +2.72%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.006'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6090'
$ws.Cells.Item(41, 5).Value = '  -3.82%  '
$ws.Cells.Item(42, 5).Value = '  +0.02%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.148'
$ws.Cells.Item(43, 5).Value = '  -3.54%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.317'
$ws.Cells.Item(44, 5).Value = '  -5.99%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.5953'
$ws.Cells.Item(45, 5).Value = '  +0.77%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.98'
$ws.Cells.Item(46, 5).Value = '  -4.54%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.717'
$ws.Cells.Item(47, 5).Value = '  -1.25%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '127.34'
$ws.Cells.Item(48, 5).Value = '  +1.04%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.206'
$ws.Cells.Item(49, 5).Value = '  +0.66%  '
$ws.Cells.Item(50, 5).Value = '  -4.81%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06708'
